$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for numeric-looking price/volume strings
# by forcing Text number format on D:E before assignment, then
# reverting the style afterwards so no stray formatting remains.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.090.54"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.875.42"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "313.54"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "0.5052"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").Value = "0.3834"
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("D9").Value = "0.08634"
$ws.Range("E9").Value = "  -7.26%  "
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").Value = "6.335"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "1.872.54"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Value = "1.004"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "7.169"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").Value = "0.00001101"
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").Value = "91.08"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").Value = "0.06636"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").Value = "18.13"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "6.103"
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("D23").Value = "28.130.94"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "2.272"
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").Value = "2.087.59"
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("D28").Value = "20.72"
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("D29").Value = "157.17"
$ws.Range("D30").Value = "125.99"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("E31").Value = "  -2.35%  "
$ws.Range("E32").Value = "  -3.64%  "
$ws.Range("D33").Value = "5.590"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("D34").Value = "3.607"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").Value = "9.690"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("D36").Value = "0.02456"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("D37").Value = "0.06576"
$ws.Range("E37").Value = "  -1.33%  "
$ws.Range("D38").Value = "0.2171"
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("E39").Value = "  -3.41%  "
$ws.Range("D40").Value = "1.245"
$ws.Range("E40").Value = "  -2.93%  "
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").Value = "0.6373"
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("D43").Value = "4.896"
$ws.Range("E43").Value = "  -2.05%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "13.27"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.5981"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "1.280"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.672"
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").Value = "1.231"
$ws.Range("E48").Value = "  +3.94%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "1.990"
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "121.28"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "80.29"
$ws.Range("E51").Value = "  +2.30%  "

# Revert D:E formatting back to the workbook default style so the
# cells carry no explicit style index, matching the original layout.
$ws.Range("D2:E51").Style = "Normal"
